# Doctor Search.xlsx — fix the "Hospital" column header.
# The shared string "Hospital " (with a trailing space) is renamed to
# "Hospital" (no trailing space), and the active selection moves to E1
# (the header cell that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Hospital"

# Move the selection to the edited header cell, matching the author's
# final cursor position recorded in the saved file.
$ws.Range("E1").Select()
